$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 96

# Columns A-D in this sheet store plain text (date/time/weekday/week-as-text
# strings), even though some of the literals ("2024-01-27", "03") look like
# a date or a number. Assigning such a literal straight to .Value makes
# Excel auto-convert it (date serial / numeric 3, dropping the leading
# zero). Route the literal through a scratch cell's text-formula result and
# Copy + PasteSpecial(values only) into the destination instead - the
# pasted value keeps its literal string identity without ever touching the
# destination cell's number format, so no extra style is introduced either.
$scratch = $ws.Cells.Item(200, 1)

function Set-TextValue($targetCell, $text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $targetCell.PasteSpecial(-4163)
}

Set-TextValue $ws.Cells.Item($row, 1) "2024-01-27"
Set-TextValue $ws.Cells.Item($row, 2) "18:17:15"
Set-TextValue $ws.Cells.Item($row, 3) "Saturday"
Set-TextValue $ws.Cells.Item($row, 4) "03"

$scratch.Clear()

$ws.Cells.Item($row, 5).Value = 137262
$ws.Cells.Item($row, 6).Value = 141691
$ws.Cells.Item($row, 7).Value = 171515
$ws.Cells.Item($row, 8).Value = 149171
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 121866
$ws.Cells.Item($row, 11).Value = 223879
$ws.Cells.Item($row, 12).Value = 257127
$ws.Cells.Item($row, 13).Value = 185335
$ws.Cells.Item($row, 14).Value = 110006
$ws.Cells.Item($row, 15).Value = 41420
$ws.Cells.Item($row, 16).Value = 30834
$ws.Cells.Item($row, 17).Value = 73613
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42392
$ws.Cells.Item($row, 20).Value = -1
